$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9806088209152222
$ws.Range("B1").Value = 1.243554830551147
$ws.Range("C1").Value = 1.068673133850098
$ws.Range("D1").Value = 1.104182839393616
$ws.Range("E1").Value = 1.235699772834778
